# 13.1.3.xlsx update
# - A1 label updated from the old "1.5.4 ..." indicator text to "13.1.3 ..."
# - D4 changes from a text-typed "484" into the literal number 484
# - Five new year columns (2020-2023 plus the existing 2019 series extended)
#   are added in columns E:H for rows 3-6, copying formatting from column D
#   and filling in the reported values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend formatting of column D into E:H for the data rows -------------
$ws.Range("D3").Copy($ws.Range("E3:H3"))
$ws.Range("D4").Copy($ws.Range("E4:H4"))
$ws.Range("D5").Copy($ws.Range("E5:H5"))
$ws.Range("D6").Copy($ws.Range("E6:H6"))

# --- row 1: corrected indicator number in the Kyrgyz label ----------------
$ws.Range("A1").Value = "13.1.3 Кырсыктардын кооптуулугун азайтуунун улуттук стратегияларына ылайык, кырсыктардын кооптуулугун азайтуунун жергиликтүү стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын үлүшү"

# --- row 3: year headers ----------------------------------------------------
$ws.Range("D3").Value = 2019
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# --- row 4: number of local governments (now a literal number, not text) --
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# --- row 5: proportion (%) -------------------------------------------------
$ws.Range("D5").Value = 10.12397
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# --- row 6: count of local governments adopting DRR strategies ------------
$ws.Range("D6").Value = 49
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169
